$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("liquids")
$ws.Activate()

# Capture the current (pre-edit) pg/vg header labels and per-row values before
# they get overwritten, then shift pg -> D, vg -> E and insert a new
# "volume" column at C.
$pgHeader = $ws.Range("C1").Value()
$vgHeader = $ws.Range("D1").Value()
$descHeader = $ws.Range("E1").Value()

$rows = 2..12
$pgValues = @{}
$vgValues = @{}
foreach ($r in $rows) {
    $pgValues[$r] = $ws.Cells.Item($r, 3).Value()
    $vgValues[$r] = $ws.Cells.Item($r, 4).Value()
}

# Volume (new column) values per row, in sheet order rows 2..12.
$volumeByRow = @{
    2 = 20
    3 = 20
    4 = 15
    5 = 10
    6 = 10
    7 = 10
    8 = 10
    9 = 20
    10 = 20
    11 = 20
    12 = 20
}

# Header row: A/B stay put, C becomes "volume", D/E take the old pg/vg
# headers, F takes the old description header (carrying over the bold
# "Nadpis 2" heading style used across the header row).
$ws.Range("C1").Value = "volume"
$ws.Range("D1").Value = $pgHeader
$ws.Range("E1").Value = $vgHeader
$ws.Range("F1").Value = $descHeader
$ws.Range("C1").Style = $ws.Range("B1").Style
$ws.Range("F1").Style = $ws.Range("E1").Style

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = $volumeByRow[$r]
    $ws.Cells.Item($r, 4).Value = $pgValues[$r]
    $ws.Cells.Item($r, 5).Value = $vgValues[$r]
}

# Column widths: A/B/C/D keep their existing widths; E and F get new
# explicit widths.
$ws.Columns("E").ColumnWidth = 8.16666666666667
$ws.Columns("F").ColumnWidth = 25.45

# Selection, matching the post-edit workbook.
$ws.Range("B7").Select()
